$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price figures as literal text in the source data (some use a
# European "thousands.decimal" double-dot style, some carry significant trailing
# zeros). Force each D cell we touch to Text format *before* writing so Excel
# keeps the exact literal string instead of re-interpreting it as a number.
$priceCells = @(
    "D2",
    "D3",
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D11",
    "D12",
    "D14",
    "D15",
    "D16",
    "D17",
    "D18",
    "D19",
    "D21",
    "D22",
    "D23",
    "D29",
    "D30",
    "D31",
    "D34",
    "D35",
    "D36",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D45",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '43.766.71'
$ws.Range("E2").Value = '  -0.09%  '

# Row 3
$ws.Range("D3").Value = '2.344.53'
$ws.Range("E3").Value = '  -0.45%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = '239.17'
$ws.Range("E5").Value = '  -0.99%  '

# Row 6
$ws.Range("D6").Value = '0.665'
$ws.Range("E6").Value = '  -3.72%  '

# Row 7
$ws.Range("D7").Value = '72.33'
$ws.Range("E7").Value = '  -6.65%  '

# Row 8
$ws.Range("E8").Value = '  -0.04%  '

# Row 9
$ws.Range("D9").Value = '0.591'
$ws.Range("E9").Value = '  -6.29%  '

# Row 10
$ws.Range("D10").Value = '0.0999'
$ws.Range("E10").Value = '  -2.13%  '

# Row 11
$ws.Range("D11").Value = '58.38'
$ws.Range("E11").Value = '  +1.80%  '

# Row 12
$ws.Range("D12").Value = '32.65'
$ws.Range("E12").Value = '  -3.76%  '

# Row 13
$ws.Range("E13").Value = '  -0.37%  '

# Row 14
$ws.Range("D14").Value = '7.22'
$ws.Range("E14").Value = '  -5.02%  '

# Row 15
$ws.Range("D15").Value = '2.694.92'
$ws.Range("E15").Value = '  -0.36%  '

# Row 16
$ws.Range("D16").Value = '16.05'
$ws.Range("E16").Value = '  -5.63%  '

# Row 17
$ws.Range("D17").Value = '0.899'
$ws.Range("E17").Value = '  -3.38%  '

# Row 18
$ws.Range("D18").Value = '2.343.85'
$ws.Range("E18").Value = '  -0.24%  '

# Row 19
$ws.Range("D19").Value = '43.717.58'
$ws.Range("E19").Value = '  -0.09%  '

# Row 20
$ws.Range("E20").Value = '  -1.57%  '

# Row 21
$ws.Range("D21").Value = '6.63'
$ws.Range("E21").Value = '  -0.75%  '

# Row 22
$ws.Range("D22").Value = '77.68'
$ws.Range("E22").Value = '  +0.12%  '

# Row 23
$ws.Range("D23").Value = '250.86'
$ws.Range("E23").Value = '  -2.61%  '

# Row 25
$ws.Range("E25").Value = '  +2.06%  '

# Row 26
$ws.Range("E26").Value = '  +2.28%  '

# Row 27
$ws.Range("E27").Value = '  -1.80%  '

# Row 28
$ws.Range("E28").Value = '  -6.10%  '

# Row 29
$ws.Range("D29").Value = '2.27'
$ws.Range("E29").Value = '  -1.40%  '

# Row 30
$ws.Range("D30").Value = '176.83'
$ws.Range("E30").Value = '  +1.07%  '

# Row 31
$ws.Range("D31").Value = '22.18'
$ws.Range("E31").Value = '  -4.00%  '

# Row 32
$ws.Range("E32").Value = '  -2.37%  '

# Row 33
$ws.Range("E33").Value = '  -2.28%  '

# Row 34
$ws.Range("D34").Value = '0.0752'
$ws.Range("E34").Value = '  -1.90%  '

# Row 35
$ws.Range("D35").Value = '5.06'
$ws.Range("E35").Value = '  -5.48%  '

# Row 36
$ws.Range("D36").Value = '5.31'
$ws.Range("E36").Value = '  -1.67%  '

# Row 37
$ws.Range("E37").Value = '  -1.13%  '

# Row 38
$ws.Range("B38").Value = 'THORChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D38").Value = '6.38'
$ws.Range("E38").Value = '  -0.95%  '

# Row 39
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = '2.35'
$ws.Range("E39").Value = '  -3.29%  '

# Row 40
$ws.Range("B40").Value = 'FTXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D40").Value = '5.66'
$ws.Range("E40").Value = '  +24.69%  '

# Row 41
$ws.Range("D41").Value = '0.0270'
$ws.Range("E41").Value = '  -2.99%  '

# Row 42
$ws.Range("D42").Value = '65.81'
$ws.Range("E42").Value = '  +15.18%  '

# Row 43
$ws.Range("D43").Value = '9.20'
$ws.Range("E43").Value = '  +0.93%  '

# Row 44
$ws.Range("E44").Value = '  -2.48%  '

# Row 45
$ws.Range("D45").Value = '18.86'
$ws.Range("E45").Value = '  -2.46%  '

# Row 46
$ws.Range("E46").Value = '  -5.04%  '

# Row 47
$ws.Range("E47").Value = '  +0.10%  '

# Row 48
$ws.Range("E48").Value = '  -3.11%  '

# Row 49
$ws.Range("D49").Value = '2.41'
$ws.Range("E49").Value = '  -5.02%  '

# Row 50
$ws.Range("D50").Value = '1.15'
$ws.Range("E50").Value = '  -3.76%  '

# Row 51
$ws.Range("B51").Value = 'HuobiToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D51").Value = '2.93'
$ws.Range("E51").Value = '  +3.26%  '
